# Update "想去人数" (interested-count) figures that changed between scrapes.
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the same
# events in the same rows, so both need to be updated. Note row 17 ends up
# with slightly different values on the two sheets (1889 vs 1890), matching
# the source data exactly.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1) updates
$sheetExhibit.Range("F2").Value  = 135
$sheetExhibit.Range("F4").Value  = 11980
$sheetExhibit.Range("F5").Value  = 1255
$sheetExhibit.Range("F11").Value = 422
$sheetExhibit.Range("F17").Value = 1889
$sheetExhibit.Range("F19").Value = 923

# 全部类型 (sheet4) updates
$sheetAll.Range("F2").Value  = 135
$sheetAll.Range("F4").Value  = 11980
$sheetAll.Range("F5").Value  = 1255
$sheetAll.Range("F11").Value = 422
$sheetAll.Range("F17").Value = 1890
$sheetAll.Range("F19").Value = 923
